$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1:AF1, copying the style from the existing header (AC1)
# so the new cells share the same cellXf as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 93   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
